# Remove the "Ver no Jupiter..." / copyright footer block (and the blank
# paragraph preceding it) that used to follow the last "Requisitos" entry
# ("LOQ4037: Química Orgânica I (Requisito fraco)"), per the site rebuild.

$d = $word.ActiveDocument

$startMarker = "LOQ4037: Qu"
$endMarker   = "Creative Commons Attribution"

# wdParagraph = 4
$wdParagraph = 4

$rFind = $d.Content
$null = $rFind.Find.Execute($startMarker, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$null = $rFind.Expand($wdParagraph)

$rEnd = $d.Content
$null = $rEnd.Find.Execute($endMarker, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$null = $rEnd.Expand($wdParagraph)

# Delete everything from the end of the "LOQ4037..." paragraph (its
# paragraph mark) through to the end of the copyright paragraph (including
# its own paragraph mark) -- this removes the blank paragraph, the
# "Ver no Jupiter..." paragraph and the "(c) 2020 ..." paragraph, while
# leaving the trailing blank paragraph and the page-break paragraph intact.
$deleteRange = $d.Range($rFind.End, $rEnd.End)
$deleteRange.Delete()
